{"js": "// Convert the two M2Doc Word \"field\" tokens (fldChar begin/instrText.../fldChar end)\n// used for `m:if ...` / `m:endif` into plain literal text runs wrapped in `{ }`\n// e.g. `{m:if self.name =    }` and `{m:endif}` \u2014 mirrors\n// TokenIteratorFieldRewriterSplit turning a live Word field into inert M2Doc\n// template syntax.\n//\n// Strategy per field:\n//   1. Locate the field (by its code text) and the paragraph that owns it.\n//   2. Delete the field (removes the fldChar begin/end + every instrText/.text\n//      run that belongs to it, in one shot).\n//   3. Re-insert the literal replacement text, as one `w:r/w:t` run per piece\n//      (so the run-splitting matches the diff), at the start of that\n//      paragraph, via insertOoxml so we control `xml:space=\"preserve\"` exactly\n//      like the reference edit.\n\nfunction runXml(text, preserveSpace) {\n  const sp = preserveSpace ? ' xml:space=\"preserve\"' : \"\";\n  const escaped = text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n  return `<w:r><w:t${sp}>${escaped}</w:t></w:r>`;\n}\n\nfunction wrapOoxml(innerRunsXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" +\n    innerRunsXml +\n    \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// Each field's replacement, expressed as [text, needsPreserveSpace] pieces \u2014\n// one piece becomes one `w:r` run, in this left-to-right order.\nconst fieldRewrites = [\n  {\n    match: \"m:if self.name =\",\n    pieces: [\n      [\"{m:if \", true],\n      [\"self.name \", true],\n      [\"=\", false],\n      [\"    \", true],\n      [\"}\", true]\n    ]\n  },\n  {\n    match: \"m:endif\",\n    pieces: [[\"{m:endif}\", true]]\n  }\n];\n\nfor (const rewrite of fieldRewrites) {\n  // 1. Find the live field whose code matches this rewrite (whitespace-\n  // insensitive compare, since field codes commonly carry extra spaces).\n  const fields = context.document.body.fields;\n  fields.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < fields.items.length; i++) {\n    fields.items[i].load(\"code\");\n  }\n  await context.sync();\n\n  let target = null;\n  for (let i = 0; i < fields.items.length; i++) {\n    const code = fields.items[i].code.replace(/\\s+/g, \" \").trim();\n    if (code === rewrite.match) {\n      target = fields.items[i];\n      break;\n    }\n  }\n  if (!target) {\n    continue; // already converted, or not present in this document\n  }\n\n  // Find the paragraph that owns this field (needed so we can re-insert the\n  // replacement text at the right place once the field itself is gone).\n  const paras = context.document.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n\n  let targetPara = null;\n  for (let i = 0; i < paras.items.length; i++) {\n    const pf = paras.items[i].fields;\n    pf.load(\"items/code\");\n    await context.sync();\n    for (let j = 0; j < pf.items.length; j++) {\n      const code = pf.items[j].code.replace(/\\s+/g, \" \").trim();\n      if (code === rewrite.match) {\n        targetPara = paras.items[i];\n        break;\n      }\n    }\n    if (targetPara) break;\n  }\n  if (!targetPara) {\n    continue;\n  }\n\n  // 2. Remove the field outright (clears fldChar begin/instrText*/fldChar end).\n  target.delete();\n  await context.sync();\n\n  // 3. Re-insert the literal text runs at the paragraph start, exactly where\n  // the field used to begin.\n  const innerRunsXml = rewrite.pieces.map(([text, preserve]) => runXml(text, preserve)).join(\"\");\n  targetPara.insertOoxml(wrapOoxml(innerRunsXml), Word.InsertLocation.start);\n  await context.sync();\n}\n", "ps1": "# Convert the two M2Doc Word \"field\" tokens (fldChar begin/instrText.../fldChar\n# end) used for `m:if ...` / `m:endif` into plain literal text runs wrapped in\n# `{ }`, e.g. `{m:if self.name =    }` and `{m:endif}` \u2014 mirrors\n# TokenIteratorFieldRewriterSplit turning a live Word field into inert M2Doc\n# template syntax.\n#\n# Strategy per field:\n#   1. Re-scan $d.Fields (collections/indices shift after every edit, so this\n#      is done fresh for every field we convert) and find the field whose\n#      Code text matches (whitespace-insensitive) what we are looking for.\n#   2. Locate the paragraph that owns it by comparing Range positions.\n#   3. Delete the field (clears fldChar begin/end + every instrText run that\n#      belongs to it, in one shot).\n#   4. Re-insert the literal replacement text pieces, one `InsertBefore` call\n#      per piece (so the run-splitting matches the reference edit), right\n#      where the field used to start.\n\n$d = $word.ActiveDocument\n\nfunction Normalize-Code($text) {\n    return ($text -replace '\\s+', ' ').Trim()\n}\n\nfunction Convert-Field($matchCode, $pieces) {\n    $d = $word.ActiveDocument\n\n    # 1. Find the live field whose code matches (collapsing whitespace runs).\n    $target = $null\n    for ($i = 1; $i -le $d.Fields.Count; $i++) {\n        $f = $d.Fields.Item($i)\n        if ((Normalize-Code $f.Code.Text) -eq $matchCode) {\n            $target = $f\n            break\n        }\n    }\n    if ($null -eq $target) {\n        return # already converted, or not present in this document\n    }\n\n    # 2. Find the paragraph that contains this field's code range.\n    $fStart = $target.Code.Start\n    $targetPara = $null\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($fStart -ge $p.Range.Start -and $fStart -lt $p.Range.End) {\n            $targetPara = $p\n            break\n        }\n    }\n    if ($null -eq $targetPara) {\n        return\n    }\n\n    # 3. Remove the field outright.\n    $target.Delete()\n\n    # 4. Re-insert the literal text pieces, in reverse order, always at the\n    # (now-collapsed) start of the paragraph, so their final left-to-right\n    # order matches $pieces.\n    for ($i = $pieces.Length - 1; $i -ge 0; $i--) {\n        $r = $targetPara.Range\n        $r.SetRange($r.Start, $r.Start)\n        $r.InsertBefore($pieces[$i])\n    }\n}\n\nConvert-Field \"m:if self.name =\" @(\"{m:if \", \"self.name \", \"=\", \"    \", \"}\")\nConvert-Field \"m:endif\" @(\"{m:endif}\")\n"}
